$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the English translation column (C) ---
$ws.Range("C1").Value = "After assigning private IP for the elastic network interface, please login to the virtual machine to configure the private IP. Click to view"
$ws.Range("C3").Value = "Unable to assign new IP since the upper limit has been reached."
$ws.Range("C6").Value = "IP Type"
$ws.Range("C8").Value = "Actions"
$ws.Range("C10").Value = "Associate EIP"
# Typing a leading apostrophe forces "text" entry in Excel; the apostrophe itself
# is consumed as a marker (quotePrefix) and is not part of the stored text, while
# the trailing apostrophe remains part of the value.
$ws.Range("C11").Value = "'Primary IP'"

# --- Layout tweaks ---
$ws.Rows.Item(1).RowHeight = 46.5
$ws.Columns.Item(3).ColumnWidth = 46.21875

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("C17").Select()
